$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width: 116.1640625 -> 118.33203125 (closest achievable via COM pixel quantization) ---
$ws.Columns("A").ColumnWidth = 117.5

# --- Rows 16-33: apply the "red / Calibri (Body)" style (same as used on rows 36-41) to column A ---
$ws.Range("A36").Copy() | Out-Null
$ws.Range("A16:A33").PasteSpecial(-4122) | Out-Null

# --- Rows 42-65: apply the "green / Calibri (Body)" style (same as used on rows 3-14) to column A ---
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A42:A65").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Fill in the previously-missing B column values (repeats of the preceding row's B value) ---
$ws.Range("B18").Value2 = $ws.Range("B17").Value2
$ws.Range("B33").Value2 = $ws.Range("B32").Value2
$ws.Range("B50").Value2 = $ws.Range("B47").Value2
$ws.Range("B65").Value2 = $ws.Range("B64").Value2

# --- Update selection / scroll position: was topLeftCell A20 + selection A44, now just selection A49 ---
$ws.Range("A1").Select() | Out-Null
$ws.Range("A49").Select() | Out-Null
